# Update the "settings" sheet: rename form_title / form_id from
# "Budgets FR" / "BudgetsFR" to "Budgets" (Kobo forms in French).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

$ws.Range("A2").Value = "Budgets"
$ws.Range("B2").Value = "Budgets"

# Move the active selection to B2 to match the saved cursor position.
$ws.Activate()
$ws.Range("B2").Select()
